# Strip the 'EMBARGO DATE' column from the "Item description" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Item description")

# Column C holds the "EMBARGO DATE" header/values; delete the whole column,
# which shifts every subsequent column one place to the left.
$ws.Columns.Item(3).Delete()
